# "Busca pela lupa falha terminado." - fix the search-by-magnifier test data.
#
# The "Pesquisa pela lupa" (search by magnifying glass) sheet held a stale
# product name ("HP PAVILION 15Z TOUCH LAPTOP") that no longer matches the
# product actually used elsewhere in the workbook. Update it to the product
# that the "Pesquisa pagina inicial" sheet already searches for/shows
# ("HP ZBook 17 G2 Mobile Workstation", typed here in upper case as the
# lupa/search field is used), picking up that sheet's cell formatting, and
# make this sheet the active one.

$wb = $excel.ActiveWorkbook

$wsInicial = $wb.Worksheets.Item("Pesquisa pagina inicial")
$wsLupa = $wb.Worksheets.Item("Pesquisa pela lupa")

# Set the new value first ...
$wsLupa.Range("A1").Value = "HP ZBOOK 17 G2 MOBILE WORKSTATION"

# ... then pick up the formatting used for the matching product string on
# the "Pesquisa pagina inicial" sheet (copy/paste formats only).
$wsInicial.Range("A2").Copy()
$wsLupa.Range("A1").PasteSpecial(-4122)

# The edit was made on (and leaves active) the "Pesquisa pela lupa" tab.
$wsLupa.Activate()
